# Updated cryptos list (price + 1h volume change columns) per the latest
# GitHub Actions refresh. Price-column ("D") values look numeric (e.g.
# "61.688,04" style strings using '.' separators) but must stay plain TEXT
# cells, matching the workbook's existing inlineStr convention - so we
# write them with a leading apostrophe (forces text entry, same as typing
# '61.688.04 into Excel) and then reset the cell Style back to "Normal" so
# no stray NumberFormat/quote-prefix formatting is left behind.
# Volume-column ("E") values already carry padding spaces around the
# percentage text, so Excel never mistakes them for numbers and they can
# be assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.688.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.32%  '
$ws.Range("D3").Value = "'3.041.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.22%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'581.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").Value = "'129.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.81%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = "'3.039.92"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.17%  '
$ws.Range("D9").Value = "'0.502"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("E10").Value = '  -2.13%  '
$ws.Range("D11").Value = "'5.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("E12").Value = '  -3.47%  '
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("D14").Value = "'33.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("D16").Value = "'3.546.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.24%  '
$ws.Range("D17").Value = "'61.711.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.26%  '
$ws.Range("D18").Value = "'3.046.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.31%  '
$ws.Range("E19").Value = '  -2.08%  '
$ws.Range("D20").Value = "'445.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.89%  '
$ws.Range("D21").Value = "'13.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.31%  '
$ws.Range("D22").Value = "'0.669"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.40%  '
$ws.Range("D23").Value = "'7.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.58%  '
$ws.Range("D24").Value = "'80.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.56%  '
$ws.Range("D25").Value = "'12.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.02%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("E28").Value = '  -4.36%  '
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("D30").Value = "'7.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.77%  '
$ws.Range("D31").Value = "'6.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.32%  '
$ws.Range("D32").Value = "'25.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.08%  '
$ws.Range("D33").Value = "'0.0966"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.04%  '
$ws.Range("E34").Value = '  -2.05%  '
$ws.Range("D35").Value = "'0.970"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.81%  '
$ws.Range("D36").Value = "'5.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.62%  '
$ws.Range("D37").Value = "'50.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").Value = "'0.0372"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.26%  '
$ws.Range("D40").Value = "'7.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.97%  '
$ws.Range("E41").Value = '  -2.08%  '
$ws.Range("E42").Value = '  -6.94%  '
$ws.Range("D43").Value = "'377.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.05%  '
$ws.Range("D44").Value = "'2.681.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.37%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = "'123.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("E47").Value = '  -4.11%  '
$ws.Range("D48").Value = "'34.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.59%  '
$ws.Range("E49").Value = '  -5.80%  '
$ws.Range("E50").Value = '  -2.45%  '
$ws.Range("D51").Value = "'23.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.70%  '